# Fix #366 - "User content is lost after two generation without edition."
#
# M2Doc's "m:usercontent" / "m:endusercontent" markers were written out as
# <w:fldSimple w:instr="..."/> elements. Re-running the generator a second
# time on such a document lost the protected user content, because a
# fldSimple cannot carry a separate field-result run. Expanding every such
# field to the classic begin/instrText/separate/end run sequence fixes it.
#
# This walks the document's fields (from last to first, so earlier inserts
# never shift the character offsets of fields still to be processed) and,
# for each field, rebuilds its containing paragraph's content as the
# expanded field-code form while preserving whatever paragraph-mark run
# formatting (character style / color) that paragraph already had.

$d = $word.ActiveDocument
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-RgbHex($colorValue) {
    $r = $colorValue -band 0xFF
    $g = ($colorValue -shr 8) -band 0xFF
    $b = ($colorValue -shr 16) -band 0xFF
    return ('{0:X2}{1:X2}{2:X2}' -f $r, $g, $b)
}

for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $field = $d.Fields.Item($i)
    $instr = $field.Code.Text.Trim()

    # Locate the paragraph that owns this field by counting paragraph marks
    # up to the field's start (more reliable here than Range.Paragraphs on
    # a collapsed point range).
    $paraIndex = $d.Range(0, $field.Code.Start).Paragraphs.Count
    $para = $d.Paragraphs.Item($paraIndex)

    # Preserve the paragraph-mark run formatting (character style + color),
    # e.g. the green "lev" style used to highlight protected user content.
    $styleName = $para.Style.NameLocal
    $pPrXml = ""
    if ($styleName -and $styleName -ne "Normal") {
        $colorHex = Get-RgbHex $para.Range.Font.Color
        $pPrXml = '<w:pPr><w:rPr><w:rStyle w:val="' + $styleName + '"/><w:color w:val="' + $colorHex + '"/></w:rPr></w:pPr>'
    }

    $fieldXml = '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' + `
        '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' + `
        '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' + `
        '<w:r><w:fldChar w:fldCharType="end"/></w:r>'

    $newParaXml = '<w:p' + $wNs + '>' + $pPrXml + $fieldXml + '</w:p>'

    $para.Range.InsertXML($newParaXml)
}
